# Apply the "Tasks" sheet update:
#  - Task 14 (row14) notes/date updated: new note + 02/24/2024
#  - New Task 15 "Tasks Factory Seeder" inserted where "Task 16: Projects Read" used to be (row15),
#    status switched from Pending -> In Development, with a note.
#  - Projects CRUD tasks renumbered/shifted down a row (Read/Create/Update/Delete).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: Task 14 Projects Page Frontend -> update note + date ---
$ws.Range("C14").Value = "Template Sent to Self in Chat fb"
# D14 stores the date as plain text (like the other date-note cells in this
# column), so force text formatting to avoid Excel auto-converting it to a
# date serial number.
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "02/24/2024"
$ws.Range("D14").Style = "Normal"

# --- Row 15: becomes the new "Task 15: Tasks Factory Seeder" row ---
$ws.Range("A15").Value = "Task 15: Tasks Factory Seeder"
$ws.Range("C15").Value = "added db seeder, to test further"
# Match the "In Development" status formatting used by rows 13/14 (bold themed font)
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("B15").Value = "In Development"

# --- Rows 16-19: renumber/shift the Projects CRUD tasks down by one ---
$ws.Range("A16").Value = "Task 16: Projects Read (CRUD)"
$ws.Range("A17").Value = "Task 17: Projects Create (CRUD)"
$ws.Range("A18").Value = "Task 18:Projects Update (CRUD)"
$ws.Range("A19").Value = "Task 19: Projects Delete (CRUD)"

# --- Selection cursor position, matches the saved workbook view ---
$ws.Range("C20").Select() | Out-Null

$excel.CutCopyMode = 0
